$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.106.83"
$ws.Range("E2").Value = "  +1.23%  "
$ws.Range("D3").Value = "2.695.84"
$ws.Range("E3").Value = "  +2.34%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'614.17"
$ws.Range("E5").Value = "  +1.73%  "
$ws.Range("D6").Value = "'158.61"
$ws.Range("E6").Value = "  +2.19%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +0.84%  "
$ws.Range("D9").Value = "'0.126"
$ws.Range("E9").Value = "  +6.34%  "
$ws.Range("E10").Value = "  +4.61%  "
$ws.Range("D11").Value = "'0.405"
$ws.Range("E11").Value = "  -0.87%  "
$ws.Range("E12").Value = "  +0.84%  "
$ws.Range("E13").Value = "  +9.86%  "
$ws.Range("D14").Value = "'30.27"
$ws.Range("E14").Value = "  +3.47%  "
$ws.Range("D15").Value = "3.177.76"
$ws.Range("E15").Value = "  +2.42%  "
$ws.Range("D16").Value = "65.950.35"
$ws.Range("E16").Value = "  +1.37%  "
$ws.Range("D17").Value = "2.686.77"
$ws.Range("E17").Value = "  +1.75%  "
$ws.Range("D18").Value = "'12.73"
$ws.Range("E18").Value = "  +1.20%  "
$ws.Range("E19").Value = "  -0.24%  "
$ws.Range("D20").Value = "'7.84"
$ws.Range("E20").Value = "  +6.80%  "
$ws.Range("D21").Value = "'358.88"
$ws.Range("E21").Value = "  -0.33%  "
$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").Value = "'71.22"
$ws.Range("E23").Value = "  +3.35%  "
$ws.Range("E24").Value = "  +18.72%  "
$ws.Range("D25").Value = "'9.95"
$ws.Range("E25").Value = "  +5.86%  "
$ws.Range("E26").Value = "  -1.46%  "
$ws.Range("E27").Value = "  +1.75%  "
$ws.Range("D28").Value = "'0.173"
$ws.Range("E28").Value = "  +4.13%  "
$ws.Range("E29").Value = "  +0.93%  "
$ws.Range("E30").Value = "  +4.29%  "
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "  +0.15%  "
$ws.Range("D32").Value = "'538.56"
$ws.Range("E32").Value = "  +2.45%  "
$ws.Range("D33").Value = "'1.80"
$ws.Range("E33").Value = "  -0.23%  "
$ws.Range("E34").Value = "  +5.31%  "
$ws.Range("E35").Value = "  -0.77%  "
$ws.Range("D36").Value = "'0.437"
$ws.Range("E36").Value = "  +2.20%  "
$ws.Range("D37").Value = "'20.86"
$ws.Range("E37").Value = "  +2.02%  "
$ws.Range("E38").Value = "  +1.65%  "
$ws.Range("E39").Value = "  -1.06%  "
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("D42").Value = "'169.23"
$ws.Range("E42").Value = "  +2.52%  "
$ws.Range("D43").Value = "'42.59"
$ws.Range("E43").Value = "  +0.50%  "
$ws.Range("E44").Value = "  +0.62%  "
$ws.Range("D45").Value = "'0.0629"
$ws.Range("E45").Value = "  +2.04%  "
$ws.Range("D46").Value = "'2.34"
$ws.Range("E46").Value = "  +4.86%  "
$ws.Range("D47").Value = "'23.81"
$ws.Range("E47").Value = "  +1.95%  "
$ws.Range("E48").Value = "  +1.01%  "
$ws.Range("D49").Value = "'0.661"
$ws.Range("E49").Value = "  +0.85%  "
$ws.Range("D50").Value = "'21.01"
$ws.Range("E50").Value = "  +7.06%  "
$ws.Range("D51").Value = "'0.0996"
$ws.Range("E51").Value = "  +1.51%  "
